# Apply the "new version with timestamp" update to the Day Sale / out-of-stock
# items report:
#   1. Insert a new item row for "ZURCAL 40MG 14 GASTRO RESISTANT TAB" right
#      before "امواس لورد".
#   2. Insert a new item row for "حبايه" right before "ريكسونا حريمي بليه".
#   3. Renumber the item (م) column, recompute the grand total, and refresh
#      the printed timestamp in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: insert the two new rows by copying an existing template row so that
# merged cells / number formats / fonts / fills come along for the ride.
# Insert from the bottom up so earlier row numbers stay valid.
# ---------------------------------------------------------------------------

# New row for "حبايه" goes directly above the current row 15 (ريكسونا حريمي بليه)
$ws.Rows.Item(15).Copy()
$ws.Rows.Item(15).Insert()

# New row for "ZURCAL ..." goes directly above the current row 12 (امواس لورد)
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(12).Insert()

# ---------------------------------------------------------------------------
# Step 2: make sure the two freshly inserted rows are merged the same way as
# every other item row (A:B, C:G, H:K, L:M, N:O).
# ---------------------------------------------------------------------------
foreach ($r in 12, 16) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

# ---------------------------------------------------------------------------
# Step 3: write out the full item table (14 rows) top to bottom so every
# value, including the ones that merely shifted down, ends up correct.
# Columns: A = م (item no.), C = name, H = الرصيد الحالي, L = حد الطلب,
#          N = السعر, P = سعر البيع, Q = عدد التعاملات
# ---------------------------------------------------------------------------
$items = @(
    @(1,  "CORASORE 150MG 20 TAB",               "1:0",    "1", "46.00",  "46.0000", "1:0", 25.5),
    @(2,  "DOLIPRANE 1 GM 15 TABS.",              "12:0",   "1", "48.00",  "48.0000", "1:0", 24.75),
    @(3,  "EREC 100MG 12 F.C. TABLETS",           "1:10",   "1", "144.00", "36.0000", "0:3", 25.5),
    @(4,  "FAWAR FRUIT 6 SACHETS",                "5:1",    "1", "24.00",  "7.9200",  "0:2", 24.75),
    @(5,  "WATER FOR INJECTION AMP. 5 ML",        "8287:0", "1", "2.00",   "2.0000",  "1:0", 25.5),
    @(6,  "ZURCAL 40MG 14 GASTRO RESISTANT TAB",  "5:0",    "1", "96.00",  "96.0000", "1:0", 25.5),
    @(7,  "امواس لورد",                            "26:0",   "0", "15.00",  "15.0000", "1:0", 24.75),
    @(8,  "ببرونه صغير الجو",                       "7:0",    "0", "20.00",  "20.0000", "1:0", 25.5),
    @(9,  "بلاستر مترسيلك 2 سم",                    "32:0",   "0", "15.00",  "15.0000", "1:0", 24.75),
    @(10, "حبايه",                                 "0:0",    "0", "3.00",   "6.0000",  "2:0", 25.5),
    @(11, "ريكسونا حريمي بليه",                     "6:0",    "0", "27.00",  "27.0000", "1:0", 25.5),
    @(12, "كالونا ",                               "0:0",    "0", "15.00",  "15.0000", "1:0", 24.75),
    @(13, "ماكينه حلاقه جليت فليكتور",               "14:0",   "0", "15.00",  "15.0000", "1:0", 25.5),
    @(14, "محلول جلوكوز 5%",                       "20:0",   "0", "27.00",  "27.0000", "1:0", 24.75)
)

$row = 7
foreach ($item in $items) {
    $ws.Range("A" + $row).Value = $item[0]
    $ws.Range("C" + $row).Value = $item[1]
    $ws.Range("H" + $row).Value = $item[2]
    $ws.Range("L" + $row).Value = $item[3]
    $ws.Range("N" + $row).Value = $item[4]
    $ws.Range("P" + $row).Value = $item[5]
    $ws.Range("Q" + $row).Value = $item[6]
    $ws.Rows.Item($row).RowHeight = $item[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Step 4: the grand total row (now row 21) and the footer row (now row 22)
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).RowHeight = 25.5
$ws.Range("P21").Value = 375.92000000000002
$ws.Rows.Item(22).RowHeight = 16.5
$ws.Range("A22").Value = "Friday, 12 September, 2025 1:49 PM"
$ws.Range("G22").Value = "1/1"
$ws.Range("K22").Value = "developed by : Abdelaziz Talaat"

Write-Host "Workbook updated"
